# Commit: "more people are working!"
#
# Two substantive changes are made to this transcript document:
#   1. The section's page size now carries an explicit orientation
#      (portrait) instead of leaving it implicit.
#   2. A footer is added to the (single) section: the classic
#      "Blank (Three Columns)" gallery footer — a borderless 3-column
#      table whose cells are left/center/right aligned, followed by a
#      trailing Footer-styled paragraph.
#
# $word / $d (ActiveDocument) are pre-seeded by the host.

$d = $word.ActiveDocument

# --- 1. Make the page orientation explicit (portrait) -------------------
$d.PageSetup.Orientation = 0   # wdOrientPortrait

# --- 2. Build the three-column footer ------------------------------------
$sec = $d.Sections.Item(1)
$footer = $sec.Footers.Item(1)

# Inserting the table consumes/extends the (until-now nonexistent)
# default footer, which mints word/footer1.xml + the footer relationship
# + footerReference automatically on save.
$tbl = $d.Tables.Add($footer.Range, 1, 3)

# Touching "Table Grid" mints the built-in TableGrid style definition
# into styles.xml (as the real gallery entry does), then we put the
# table itself back on the plain "Normal Table" style it should render
# with. Re-fetch the footer/table handles after every mutation: this
# host invalidates stale anchors once a structural edit lands.
$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Style = "Table Grid"

$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Style = "Normal Table"

# Cell 1 — left-aligned tab slot.
$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Cell(1, 1).Range.Style = "Header"

$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Cell(1, 1).Range.ParagraphFormat.Alignment = 0   # wdAlignParagraphLeft

# Cell 2 — center tab slot.
$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Cell(1, 2).Range.Style = "Header"

$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Cell(1, 2).Range.ParagraphFormat.Alignment = 1   # wdAlignParagraphCenter

# Cell 3 — right tab slot.
$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Cell(1, 3).Range.Style = "Header"

$footer2 = $sec.Footers.Item(1)
$tbl2 = $footer2.Range.Tables.Item(1)
$tbl2.Cell(1, 3).Range.ParagraphFormat.Alignment = 2   # wdAlignParagraphRight

# Trailing Footer-styled paragraph after the table (every footer ends in
# a normal paragraph, never directly in a table).
$footer2 = $sec.Footers.Item(1)
$endRng = $footer2.Range
$endRng.Collapse(0)
$endRng.InsertParagraphAfter()
